$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data rows 2-9 (A:T) with refreshed TPM-derived NATMI values
# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Tnfsf11"
$ws.Cells.Item(2,3).Value = "Tnfrsf11b"
$ws.Cells.Item(2,4).Value = "FAPs"
$ws.Cells.Item(2,5).Value = 1.0
$ws.Cells.Item(2,6).Value = 0.3333333333333333
$ws.Cells.Item(2,7).Value = 0.020131
$ws.Cells.Item(2,8).Value = 0.060393
$ws.Cells.Item(2,9).Value = 0.0058831740909272
$ws.Cells.Item(2,10).Value = 0.006076836619800507
$ws.Cells.Item(2,11).Value = 3.0
$ws.Cells.Item(2,12).Value = 1.0
$ws.Cells.Item(2,13).Value = 2.311298
$ws.Cells.Item(2,14).Value = 6.933894
$ws.Cells.Item(2,15).Value = 0.9776743782722859
$ws.Cells.Item(2,16).Value = 0.9850046587426607
$ws.Cells.Item(2,17).Value = 0.046528740038
$ws.Cells.Item(2,18).Value = 0.418758660342
$ws.Cells.Item(2,19).Value = 0.005751828571614871
$ws.Cells.Item(2,20).Value = 0.005985712380921502

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Tnfsf11"
$ws.Cells.Item(3,3).Value = "Tnfrsf11b"
$ws.Cells.Item(3,4).Value = "MuSCs"
$ws.Cells.Item(3,5).Value = 1.0
$ws.Cells.Item(3,6).Value = 0.3333333333333333
$ws.Cells.Item(3,7).Value = 0.020131
$ws.Cells.Item(3,8).Value = 0.060393
$ws.Cells.Item(3,9).Value = 0.0058831740909272
$ws.Cells.Item(3,10).Value = 0.006076836619800507
$ws.Cells.Item(3,11).Value = 1.0
$ws.Cells.Item(3,12).Value = 0.5
$ws.Cells.Item(3,13).Value = 0.0527795
$ws.Cells.Item(3,14).Value = 0.105559
$ws.Cells.Item(3,15).Value = 0.02232562172771408
$ws.Cells.Item(3,16).Value = 0.01499534125733917
$ws.Cells.Item(3,17).Value = 0.0010625041145
$ws.Cells.Item(3,18).Value = 0.006375024687
$ws.Cells.Item(3,19).Value = 0.0001313455193123289
$ws.Cells.Item(3,20).Value = 0.00009112423887900405

# Row 4
$ws.Cells.Item(4,1).Value = "FAPs"
$ws.Cells.Item(4,2).Value = "Tnfsf11"
$ws.Cells.Item(4,3).Value = "Tnfrsf11b"
$ws.Cells.Item(4,4).Value = "FAPs"
$ws.Cells.Item(4,5).Value = 3.0
$ws.Cells.Item(4,6).Value = 1.0
$ws.Cells.Item(4,7).Value = 3.059758333333333
$ws.Cells.Item(4,8).Value = 9.179275
$ws.Cells.Item(4,9).Value = 0.8941975535822989
$ws.Cells.Item(4,10).Value = 0.9236327796800838
$ws.Cells.Item(4,11).Value = 3.0
$ws.Cells.Item(4,12).Value = 1.0
$ws.Cells.Item(4,13).Value = 2.311298
$ws.Cells.Item(4,14).Value = 6.933894
$ws.Cells.Item(4,15).Value = 0.9776743782722859
$ws.Cells.Item(4,16).Value = 0.9850046587426607
$ws.Cells.Item(4,17).Value = 7.072013316316668
$ws.Cells.Item(4,18).Value = 63.64811984685001
$ws.Cells.Item(4,19).Value = 0.8742340372511731
$ws.Cells.Item(4,20).Value = 0.9097825909523162

# Row 5
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Tnfsf11"
$ws.Cells.Item(5,3).Value = "Tnfrsf11b"
$ws.Cells.Item(5,4).Value = "MuSCs"
$ws.Cells.Item(5,5).Value = 3.0
$ws.Cells.Item(5,6).Value = 1.0
$ws.Cells.Item(5,7).Value = 3.059758333333333
$ws.Cells.Item(5,8).Value = 9.179275
$ws.Cells.Item(5,9).Value = 0.8941975535822989
$ws.Cells.Item(5,10).Value = 0.9236327796800838
$ws.Cells.Item(5,11).Value = 1.0
$ws.Cells.Item(5,12).Value = 0.5
$ws.Cells.Item(5,13).Value = 0.0527795
$ws.Cells.Item(5,14).Value = 0.105559
$ws.Cells.Item(5,15).Value = 0.02232562172771408
$ws.Cells.Item(5,16).Value = 0.01499534125733917
$ws.Cells.Item(5,17).Value = 0.1614925149541667
$ws.Cells.Item(5,18).Value = 0.9689550897250001
$ws.Cells.Item(5,19).Value = 0.01996351633112575
$ws.Cells.Item(5,20).Value = 0.01385018872776762

# Row 6
$ws.Cells.Item(6,1).Value = "Inflammatory-Mac"
$ws.Cells.Item(6,2).Value = "Tnfsf11"
$ws.Cells.Item(6,3).Value = "Tnfrsf11b"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 1.0
$ws.Cells.Item(6,6).Value = 0.3333333333333333
$ws.Cells.Item(6,7).Value = 0.014756
$ws.Cells.Item(6,8).Value = 0.044268
$ws.Cells.Item(6,9).Value = 0.004312359887026068
$ws.Cells.Item(6,10).Value = 0.004454314299427563
$ws.Cells.Item(6,11).Value = 3.0
$ws.Cells.Item(6,12).Value = 1.0
$ws.Cells.Item(6,13).Value = 2.311298
$ws.Cells.Item(6,14).Value = 6.933894
$ws.Cells.Item(6,15).Value = 0.9776743782722859
$ws.Cells.Item(6,16).Value = 0.9850046587426607
$ws.Cells.Item(6,17).Value = 0.034105513288
$ws.Cells.Item(6,18).Value = 0.306949619592
$ws.Cells.Item(6,19).Value = 0.004216083771434556
$ws.Cells.Item(6,20).Value = 0.004387520336440201

# Row 7
$ws.Cells.Item(7,1).Value = "Inflammatory-Mac"
$ws.Cells.Item(7,2).Value = "Tnfsf11"
$ws.Cells.Item(7,3).Value = "Tnfrsf11b"
$ws.Cells.Item(7,4).Value = "MuSCs"
$ws.Cells.Item(7,5).Value = 1.0
$ws.Cells.Item(7,6).Value = 0.3333333333333333
$ws.Cells.Item(7,7).Value = 0.014756
$ws.Cells.Item(7,8).Value = 0.044268
$ws.Cells.Item(7,9).Value = 0.004312359887026068
$ws.Cells.Item(7,10).Value = 0.004454314299427563
$ws.Cells.Item(7,11).Value = 1.0
$ws.Cells.Item(7,12).Value = 0.5
$ws.Cells.Item(7,13).Value = 0.0527795
$ws.Cells.Item(7,14).Value = 0.105559
$ws.Cells.Item(7,15).Value = 0.02232562172771408
$ws.Cells.Item(7,16).Value = 0.01499534125733917
$ws.Cells.Item(7,17).Value = 0.0007788143020000001
$ws.Cells.Item(7,18).Value = 0.004672885812
$ws.Cells.Item(7,19).Value = 0.00009627611559151184
$ws.Cells.Item(7,20).Value = 0.00006679396298736197

# Row 8
$ws.Cells.Item(8,1).Value = "MuSCs"
$ws.Cells.Item(8,2).Value = "Tnfsf11"
$ws.Cells.Item(8,3).Value = "Tnfrsf11b"
$ws.Cells.Item(8,4).Value = "FAPs"
$ws.Cells.Item(8,5).Value = 2.0
$ws.Cells.Item(8,6).Value = 1.0
$ws.Cells.Item(8,7).Value = 0.327147
$ws.Cells.Item(8,8).Value = 0.654294
$ws.Cells.Item(8,9).Value = 0.0956069124397477
$ws.Cells.Item(8,10).Value = 0.06583606940068805
$ws.Cells.Item(8,11).Value = 3.0
$ws.Cells.Item(8,12).Value = 1.0
$ws.Cells.Item(8,13).Value = 2.311298
$ws.Cells.Item(8,14).Value = 6.933894
$ws.Cells.Item(8,15).Value = 0.9776743782722859
$ws.Cells.Item(8,16).Value = 0.9850046587426607
$ws.Cells.Item(8,17).Value = 0.7561342068060002
$ws.Cells.Item(8,18).Value = 4.536805240836
$ws.Cells.Item(8,19).Value = 0.0934724286780632
$ws.Cells.Item(8,20).Value = 0.06484883507298286

# Row 9
$ws.Cells.Item(9,1).Value = "MuSCs"
$ws.Cells.Item(9,2).Value = "Tnfsf11"
$ws.Cells.Item(9,3).Value = "Tnfrsf11b"
$ws.Cells.Item(9,4).Value = "MuSCs"
$ws.Cells.Item(9,5).Value = 2.0
$ws.Cells.Item(9,6).Value = 1.0
$ws.Cells.Item(9,7).Value = 0.327147
$ws.Cells.Item(9,8).Value = 0.654294
$ws.Cells.Item(9,9).Value = 0.0956069124397477
$ws.Cells.Item(9,10).Value = 0.06583606940068805
$ws.Cells.Item(9,11).Value = 1.0
$ws.Cells.Item(9,12).Value = 0.5
$ws.Cells.Item(9,13).Value = 0.0527795
$ws.Cells.Item(9,14).Value = 0.105559
$ws.Cells.Item(9,15).Value = 0.02232562172771408
$ws.Cells.Item(9,16).Value = 0.01499534125733917
$ws.Cells.Item(9,17).Value = 0.0172666550865
$ws.Cells.Item(9,18).Value = 0.069066620346
$ws.Cells.Item(9,19).Value = 0.002134483761684489
$ws.Cells.Item(9,20).Value = 0.0009872343277051825

# Old data had 9 data rows (10 incl. header); new data only has 8 data rows, so remove the last row
$ws.Rows.Item(10).Delete()
